# Update "想去人数" (F column) figures for both the "展览" sheet and the
# mirrored "全部类型" sheet (whose rows are offset by +1 relative to 展览).
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# row => new F value, for the 展览 sheet
$updates1 = @{
    7  = 669
    9  = 16
    10 = 5
    13 = 1566
    14 = 5358
    16 = 265
    17 = 230
    18 = 40
    19 = 15
    21 = 4464
    22 = 225
    27 = 63
}

foreach ($row in $updates1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updates1[$row]
}

# row => new F value, for the 全部类型 sheet (same data, rows shifted by +1)
$updates4 = @{
    8  = 669
    10 = 16
    11 = 5
    14 = 1566
    15 = 5358
    17 = 265
    18 = 230
    19 = 40
    20 = 15
    22 = 4464
    23 = 225
    28 = 63
}

foreach ($row in $updates4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updates4[$row]
}
